# Add 2022-Q3 data: insert a new quarter sheet after "总计" (the totals/summary
# sheet), populate it with the new quarter's fund-holding data, and record the
# new quarter in the "总计" summary sheet (pushing the existing rows down by one).

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q1    = $wb.Worksheets.Item("2022-Q1")   # template: same look & feel as the new quarter sheet

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet right after "总计" by duplicating the
#    "2022-Q1" sheet (keeps identical formatting: styles, page setup, etc.)
#    then renaming it and trimming/overwriting its data.
#    Tab order becomes: 总计, 2022-Q3, 2022-Q1, 2021-Q4, 2021-Q2.
# ---------------------------------------------------------------------------
$q1.Copy($null, $total)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template had 3 data rows (2022-Q1 had 3 funds); 2022-Q3 only has 2, so
# drop the extra row.
$q3.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# 2) Overwrite the 2022-Q3 fund holding data (2 rows).
# ---------------------------------------------------------------------------

# Row 2: 164811 工银瑞信中证京津冀协同发展主题指数（LOF）A
$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "'164811"
$q3.Cells.Item(2, 3).Value = "工银瑞信中证京津冀协同发展主题指数（LOF）A"
$q3.Cells.Item(2, 4).Value = "'0.12"
$q3.Cells.Item(2, 5).Value = "'93.09"
$q3.Cells.Item(2, 6).Value = "'3.78"
$q3.Cells.Item(2, 7).Value = "'0.0045"
$q3.Cells.Item(2, 8).Value = 1

# Row 3: 164825 工银瑞信中证京津冀协同发展主题指数（LOF）C
$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "'164825"
$q3.Cells.Item(3, 3).Value = "工银瑞信中证京津冀协同发展主题指数（LOF）C"
$q3.Cells.Item(3, 4).Value = "'0.03"
$q3.Cells.Item(3, 5).Value = "'93.09"
$q3.Cells.Item(3, 6).Value = "'3.78"
$q3.Cells.Item(3, 7).Value = "'0.0011"
$q3.Cells.Item(3, 8).Value = 1

# ---------------------------------------------------------------------------
# 3) Update the "总计" (totals) sheet: insert a new row for 2022-Q3 at the top
#    of the data (row 2), pushing the other quarters down by one row.
# ---------------------------------------------------------------------------

# Prepare new row 5 with the same style as the other index cells in column A.
$total.Cells.Item(2, 1).Copy($total.Cells.Item(5, 1))

# Shift rows 2,3,4 -> 3,4,5 (labels/counts/amounts); column A keeps 0,1,2,3.
$total.Cells.Item(5, 2).Value = "'2021-Q2"
$total.Cells.Item(5, 3).Value = 3
$total.Cells.Item(5, 4).Value = "'0.01"
$total.Cells.Item(5, 1).Value = 3

$total.Cells.Item(4, 2).Value = "'2021-Q4"
$total.Cells.Item(4, 3).Value = 3
$total.Cells.Item(4, 4).Value = "'0.01"

$total.Cells.Item(3, 2).Value = "'2022-Q1"
$total.Cells.Item(3, 3).Value = 3
$total.Cells.Item(3, 4).Value = "'0.01"

# New row 2: the 2022-Q3 summary entry.
$total.Cells.Item(2, 2).Value = "'2022-Q3"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = "'0.01"
